# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) text cells are updated in place; a handful of
# rows also swapped ranking position, so Coin (B) / Link (C) are rewritten too.
# A few Price values (D21, D31, D45) round-trip through Excel's numeric
# auto-detection and would lose a trailing zero (e.g. "6.60" -> 6.6), so those
# cells are pre-formatted as Text to preserve the exact original string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.520.09"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "2.289.38"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "320.44"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "108.05"
$ws.Range("E6").Value = "  +8.17%  "
$ws.Range("D7").Value = "0.593"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "0.578"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").Value = "39.38"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "7.99"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "0.893"
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").Value = "2.635.45"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "14.77"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "2.287.09"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "44.365.88"
$ws.Range("E18").Value = "  +3.51%  "
$ws.Range("D19").Value = "14.24"
$ws.Range("E19").Value = "  -8.00%  "
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +4.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("D22").Value = "66.83"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "3.24"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "240.12"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "10.33"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").Value = "39.51"
$ws.Range("E28").Value = "  +15.14%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "6.61"
$ws.Range("E30").Value = "  +3.81%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.80"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "163.68"
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "3.31"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("D36").Value = "2.09"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +12.29%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.05"
$ws.Range("E38").Value = "  +6.54%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "4.53"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").Value = "0.0331"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "15.46"
$ws.Range("E42").Value = "  +24.50%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "1.779.47"
$ws.Range("E44").Value = "  -7.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.210"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "5.53"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "86.66"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").Value = "76.17"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "60.45"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "8.83"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").Value = "1.73"
$ws.Range("E51").Value = "  +7.03%  "
